$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.991.25"
$ws.Range("E2").Value = "  +2.44%  "

$ws.Range("D3").Value = "2.054.42"
$ws.Range("E3").Value = "  +1.91%  "

$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").Value = "230.02"
$ws.Range("E5").Value = "  +1.68%  "

$ws.Range("D6").Value = "0.616"

$ws.Range("D7").Value = "58.21"
$ws.Range("E7").Value = "  +6.10%  "

$ws.Range("E8").Value = "  -0.13%  "

$ws.Range("E9").Value = "  +2.95%  "

$ws.Range("E10").Value = "  +2.93%  "

$ws.Range("E11").Value = "  +0.70%  "

$ws.Range("D12").Value = "2.357.67"
$ws.Range("E12").Value = "  +1.86%  "

$ws.Range("D13").Value = "14.60"
$ws.Range("E13").Value = "  +3.56%  "

$ws.Range("D14").Value = "20.64"
$ws.Range("E14").Value = "  +2.20%  "

$ws.Range("D15").Value = "0.753"
$ws.Range("E15").Value = "  +1.68%  "

$ws.Range("E16").Value = "  +2.82%  "

$ws.Range("D17").Value = "2.051.56"
$ws.Range("E17").Value = "  +1.18%  "

$ws.Range("D18").Value = "37.914.56"
$ws.Range("E18").Value = "  +2.44%  "

$ws.Range("E19").Value = "  -0.49%  "

$ws.Range("D20").Value = "69.75"
$ws.Range("E20").Value = "  +1.29%  "

$ws.Range("D21").Value = "0.0`u{2083}0830"
$ws.Range("E21").Value = "  +1.77%  "

$ws.Range("D22").Value = "224.62"
$ws.Range("E22").Value = "  +0.69%  "

$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("E24").Value = "  +0.85%  "

$ws.Range("E25").Value = "  +2.83%  "

$ws.Range("D26").Value = "9.31"
$ws.Range("E26").Value = "  +1.79%  "

$ws.Range("D27").Value = "166.37"
$ws.Range("E27").Value = "  +0.09%  "

$ws.Range("D28").Value = "0.133"
$ws.Range("E28").Value = "  +6.77%  "

$ws.Range("D29").Value = "19.03"
$ws.Range("E29").Value = "  +1.86%  "

$ws.Range("E30").Value = "  +0.14%  "

$ws.Range("E31").Value = "  +1.90%  "

$ws.Range("E32").Value = "  +0.65%  "

$ws.Range("E33").Value = "  +4.02%  "

$ws.Range("E34").Value = "  +0.13%  "

$ws.Range("E35").Value = "  +8.00%  "

$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("D37").Value = "5.97"
$ws.Range("E37").Value = "  +12.45%  "

$ws.Range("D38").Value = "3.31"
$ws.Range("E38").Value = "  +4.74%  "

$ws.Range("E39").Value = "  -0.21%  "

$ws.Range("D40").Value = "98.46"
$ws.Range("E40").Value = "  +3.87%  "

$ws.Range("E41").Value = "  +1.50%  "

$ws.Range("D42").Value = "1.481.27"
$ws.Range("E42").Value = "  +0.23%  "

$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").Value = "2.86"
$ws.Range("E43").Value = "  +3.54%  "

$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").Value = "0.0936"
$ws.Range("E44").Value = "  +2.52%  "

$ws.Range("D45").Value = "16.63"
$ws.Range("E45").Value = "  +2.24%  "

$ws.Range("D47").Value = "4.12"
$ws.Range("E47").Value = "  +17.33%  "

$ws.Range("E48").Value = "  +0.98%  "

$ws.Range("E49").Value = "  +1.50%  "

$ws.Range("E50").Value = "  -1.63%  "

$ws.Range("D51").Value = "2.247.03"
$ws.Range("E51").Value = "  +2.13%  "
